$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2.03
$ws.Range("R3").Value = 1.83
